$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: replace "John Paul" with "Lewis Shella", change Gender to FEMALE ---
# (email cell keeps the already-existing "johnpaul@gmail.com" shared string)
$ws.Range("A3").Value = "Lewis Shella"

# --- Row 4: replace "Micheal Petterson" / his email with "Petterson Pat" / new email ---
$ws.Range("A4").Value = "Petterson Pat"
$ws.Range("B4").Value = "petterson@gmail.com"

# --- Row 3 Gender (introduced after row4 strings so shared-string order matches) ---
$ws.Range("D3").Value = "FEMALE"

# --- New row 5 data (Abram Joe, developer in Yaounde) ---
$ws.Range("A5").Value = "Abram Joe"
$ws.Range("B5").Value = "joe@gmail.com"
$ws.Range("E5").Value = "Yaounde"
$ws.Range("F5").Value = "Developer"

# --- Reuse of already-existing shared strings ---
$ws.Range("B3").Value = "johnpaul@gmail.com"
$ws.Range("E4").Value = "Buea Town"
$ws.Range("F4").Value = "Engineer"
$ws.Range("C5").Value = 678934023
$ws.Range("D5").Value = "MALE"

# --- Formatting to mirror row 2/3's telephone-column style and row4's layout ---
$ws.Range("C2").Copy()
$ws.Range("C5").PasteSpecial(-4122)

$ws.Rows.Item(1).RowHeight = 13.2
$ws.Rows.Item(2).RowHeight = 13.2
$ws.Rows.Item(3).RowHeight = 13.2
$ws.Rows.Item(5).RowHeight = 15.75

# --- Hyperlink for the new email address ---
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:joe@gmail.com")
$ws.Range("B5").Style = $ws.Range("B4").Style

# --- Final selection matches the authored workbook (cell F5 active) ---
[void]$ws.Range("F5").Select()
